$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.221.34'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.640.63'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '2.639.57'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('E10').Value = '  +8.83%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.25'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000193'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.32%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').Value = '3.121.15'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').Value = '68.203.74'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '2.625.92'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '363.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.70%  '
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('E22').Value = '  +4.08%  '
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  +3.77%  '
$ws.Range('D29').Value = '2.773.29'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '563.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('E35').Value = '  +2.69%  '
$ws.Range('E37').Value = '  +4.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.30'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('E40').Value = '  +2.31%  '
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.21%  '
$ws.Range('E43').Value = '  +4.11%  '
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '157.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.48%  '
$ws.Range('E49').Value = '  +2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('E51').Value = '  +1.52%  '
